# "make robust to only one entry"
# The Creditor ID / IBAN example value was too short to be realistic and
# needs to be updated to a full-length IBAN. Since the new value is much
# longer than the old placeholder, wrap text is enabled on those cells so
# the column keeps displaying correctly with a single (or differing)
# number of entries.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$newIban = "AT611234567890123456 "

# Update both sample rows (Creditor ID / IBAN column) to the new value.
$ws.Range("E2").Value = $newIban
$ws.Range("E3").Value = $newIban

# Enable wrap text on the IBAN cells so the longer value is fully visible.
$ws.Range("E2:E3").WrapText = $true

# Move/restore the active selection to E3.
$null = $ws.Range("E3").Select()
